$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.09"
$ws.Range("E2").Value = "'-1.16%"
$ws.Range("D3").Value = "'44.38"
$ws.Range("E3").Value = "'-0.67%"
$ws.Range("D4").Value = "'5.359"
$ws.Range("E4").Value = "'-3.73%"
$ws.Range("D5").Value = "'0.08371"
$ws.Range("E5").Value = "'1.05%"
$ws.Range("D6").Value = "'1.939"
$ws.Range("E6").Value = "'-4.96%"
$ws.Range("D7").Value = "'0.9731"
$ws.Range("E7").Value = "'-0.36%"
$ws.Range("D8").Value = "'2.510"
$ws.Range("E8").Value = "'-5.06%"
$ws.Range("D9").Value = "'0.1134"
$ws.Range("E9").Value = "'1.08%"
$ws.Range("D10").Value = "'0.1902"
$ws.Range("E10").Value = "'-0.66%"
$ws.Range("D11").Value = "'0.09712"
$ws.Range("E11").Value = "'-3.47%"
$ws.Range("D12").Value = "'0.04606"
$ws.Range("E12").Value = "'-1.62%"
$ws.Range("D13").Value = "'0.1062"
$ws.Range("E13").Value = "'0.40%"
$ws.Range("D14").Value = "'0.001297"
$ws.Range("E14").Value = "'2.06%"
$ws.Range("D15").Value = "'0.005860"
$ws.Range("E15").Value = "'-3.78%"
$ws.Range("D16").Value = "'3.361"
$ws.Range("E16").Value = "'-0.03%"
$ws.Range("D17").Value = "'4.419"
$ws.Range("D18").Value = "'0.3359"
$ws.Range("D19").Value = "'8.302"
$ws.Range("E19").Value = "'-19.39%"
$ws.Range("E20").Value = "'0.14%"
$ws.Range("E21").Value = "'6.51%"
$ws.Range("E22").Value = "'1.60%"
$ws.Range("D23").Value = "'0.001239"
$ws.Range("E23").Value = "'-4.85%"
$ws.Range("D24").Value = "'0.004428"
$ws.Range("E24").Value = "'0.80%"
$ws.Range("E25").Value = "'1.58%"
$ws.Range("D26").Value = "'0.0002980"
$ws.Range("E26").Value = "'-20.34%"
$ws.Range("D38").Value = "'0.02715"
$ws.Range("E38").Value = "'-2.87%"
$ws.Range("D39").Value = "'0.05635"
$ws.Range("E39").Value = "'-1.91%"
$ws.Range("D40").Value = "'0.007833"
$ws.Range("E40").Value = "'2.67%"
$ws.Range("D41").Value = "'0.1416"
$ws.Range("E41").Value = "'-0.40%"
$ws.Range("D42").Value = "'0.007310"
$ws.Range("E42").Value = "'-3.40%"
$ws.Range("D43").Value = "'0.002041"
$ws.Range("E43").Value = "'3.39%"
$ws.Range("D44").Value = "'0.008699"
$ws.Range("E44").Value = "'4.61%"
$ws.Range("D45").Value = "'0.3509"
$ws.Range("D46").Value = "'0.00006918"
$ws.Range("E46").Value = "'-1.75%"
$ws.Range("E47").Value = "'0.01%"
$ws.Range("D48").Value = "'0.003487"
$ws.Range("E48").Value = "'-2.88%"
$ws.Range("E49").Value = "'39.91%"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("E51").Value = "'0.01%"
